$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.566.90"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.015.82"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'378.58"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'103.05"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").Value = "'0.545"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "3.492.66"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "'18.48"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "3.014.87"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "'10.43"
$ws.Range("E18").Value = "  -14.53%  "
$ws.Range("D19").Value = "51.562.99"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "'12.46"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'69.85"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'267.53"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").Value = "'8.23"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'7.52"
$ws.Range("E27").Value = "  +6.13%  "
$ws.Range("E28").Value = "  +5.22%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'26.16"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").Value = "'0.0457"
$ws.Range("E33").Value = "  +5.43%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'34.08"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.71"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("D39").Value = "'17.20"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("D40").Value = "'0.285"
$ws.Range("E40").Value = "  +9.30%  "
$ws.Range("E41").Value = "  +4.95%  "
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'126.79"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.74"
$ws.Range("E45").Value = "  +8.87%  "
$ws.Range("D46").Value = "'21.59"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "2.030.30"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "3.315.46"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "'0.0320"
$ws.Range("E51").Value = "  +1.54%  "
